$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) "Data de entrega: 13/07" -> "Data de entrega: 15/07", expressed as
#    three separate runs ("1" / "5" / "/07") each carrying the original
#    sz=32 / szCs=32 run formatting, plus the _GoBack bookmark right
#    after it (end of that paragraph).
# ----------------------------------------------------------------------

$find1 = $d.Content
$found1 = $find1.Find.Execute("13/07", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "15/07", 2)

# Split "15/07" into 3 runs ("1" | "5" | "/07") without losing formatting:
# toggling Bold on then back off on the middle character forces the
# engine to re-serialize the text as separate runs while every piece
# keeps the same resulting (non-bold) rPr.
$midChar = $d.Range($find1.Start + 1, $find1.Start + 2)   # the "5"
$midChar.Font.Bold = $true
$midChar.Font.Bold = $false

# Move the (single, document-wide) "_GoBack" bookmark to sit right after
# "/07", i.e. at the end of this paragraph's text. A bookmark created
# directly on a collapsed range landing exactly at that boundary is
# mis-anchored by this host, so: insert a throw-away character there,
# wrap the bookmark around it (non-collapsed range == reliable anchor),
# then delete the character again, leaving the bookmark collapsed in
# the right spot.
$dateParaEnd = $find1.Start + 5   # length of "15/07"
$tail = $d.Range($dateParaEnd, $dateParaEnd)
$tail.InsertAfter("X")
$wrap = $d.Range($dateParaEnd, $dateParaEnd + 1)
$d.Bookmarks.Add("_GoBack", $wrap)
$placeholder = $d.Range($dateParaEnd, $dateParaEnd + 1)
$placeholder.Text = ""

# ----------------------------------------------------------------------
# 2) Remove the bookmark that used to sit between "...OBS Studio" and
#    ", que grava a tela do computador. " (already gone - moved above)
#    and fuse those two runs of identical (bold, sz=32/szCs=32)
#    formatting back into a single run.
# ----------------------------------------------------------------------

$findObs = $d.Content
$foundObs = $findObs.Find.Execute("OBS Studio", $true, $false, $false, $false, `
                                   $false, $true, 1, $false, "", 0)
$obsStart = $findObs.Start

$findTail = $d.Content
$foundTail = $findTail.Find.Execute(", que grava a tela do computador. ", $true, `
                                     $false, $false, $false, $false, $true, 1, `
                                     $false, "", 0)
$tailEnd = $findTail.End

$mergeRange = $d.Range($obsStart, $tailEnd)
$combinedText = $mergeRange.Text

# Force the two adjacent runs to collapse into one by rewriting the
# span through a different placeholder value and then back to the
# original text.
$mergeRange.Text = "."
$mergeRange.Text = $combinedText
